# Update the "想去人数" (column F) values on each sheet to match the
# newly scraped data output (commit: "Update gh-pages to output generated
# at 456a3b4").

$wb = $excel.ActiveWorkbook

# Map of sheet name -> list of (row, newValue) pairs for column F.
$updates = @{
    "展览"     = @{ 3 = 1751; 5 = 506; 6 = 1298; 7 = 382; 9 = 917; 10 = 729;
                     16 = 3028; 18 = 526; 20 = 14; 22 = 246; 24 = 5429;
                     25 = 596; 29 = 362; 30 = 1146 }
    "演出"     = @{ 3 = 1160; 4 = 8; 13 = 621; 25 = 4000; 32 = 176 }
    "本地生活" = @{ 5 = 2516; 6 = 1080; 9 = 1373; 10 = 382 }
    "全部类型" = @{ 5 = 2516; 7 = 1751; 9 = 1080; 10 = 1373; 11 = 382;
                     14 = 506; 15 = 1298; 16 = 382; 17 = 917; 18 = 729;
                     19 = 1160; 23 = 3028; 26 = 14; 29 = 246; 31 = 5429;
                     32 = 596; 34 = 621; 37 = 362; 47 = 176 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowsForSheet = $updates[$sheetName]
    foreach ($row in $rowsForSheet.Keys) {
        $newValue = $rowsForSheet[$row]
        $ws.Range("F$row").Value = $newValue
    }
}
